$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.863.03'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.94%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.292.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.79%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.22%  '

$ws.Range("E6").Value = '  +0.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.88'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.70%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  +5.28%  '

$ws.Range("E10").Value = '  +4.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.71'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +15.30%  '

$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.631.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.818'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.289.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.798.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0954'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.85%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.54%  '

$ws.Range("E24").Value = '  +9.80%  '

$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.86%  '

$ws.Range("E27").Value = '  +1.53%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.139'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.05%  '

$ws.Range("E30").Value = '  +3.30%  '

$ws.Range("E31").Value = '  +2.45%  '

$ws.Range("E32").Value = '  +4.10%  '

$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0701'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.34%  '

$ws.Range("E35").Value = '  +1.50%  '

$ws.Range("E36").Value = '  +1.39%  '

$ws.Range("E37").Value = '  +2.40%  '

$ws.Range("E38").Value = '  +0.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("E40").Value = '  +3.88%  '

$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +27.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000220'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.73%  '

$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0970'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.52'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.10'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.489.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.35'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.60%  '
